# Generate Report for Handoff
#
# - Status moves from "Handed back: in sync with en-US" to "Ready for handoff"
#   on every sheet that tracks it (Overview!E2/F2, zh-cn!C2, de-de!C2).
# - The "Latest Handoff Datetime" timestamps advance a few seconds to the
#   freshly generated report time (Overview!G2, zh-cn!H2, de-de!H2).
# - The zh-cn/de-de "Status" column (and the Overview zh-cn/de-de columns)
#   are narrower now that the status text is shorter.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

# ---- Overview sheet ------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("G2").Value = "2016-08-26 06:57:57"

$wsOverview.Columns.Item(5).ColumnWidth = 17.2159881591797
$wsOverview.Columns.Item(6).ColumnWidth = 17.2159881591797

# ---- zh-cn sheet -----------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("H2").Value = "2016-08-26 06:57:52"

$wsZhCn.Columns.Item(3).ColumnWidth = 17.2159881591797

# ---- de-de sheet -----------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("H2").Value = "2016-08-26 06:57:57"

$wsDeDe.Columns.Item(3).ColumnWidth = 17.2159881591797
